# Weekly update: insert a new price record (row 199) for "Papa" at the
# Terminal Hortofrutícola Agro Chillán market, shifting the existing rows
# 199-230 down to 200-231, and refresh row 198 with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 199; rows 199-230 shift to 200-231.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the record that used to live in
# row 198 before this week's figures overwrote it.
$ws.Cells.Item(199, 1).Value = 7
$ws.Cells.Item(199, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(199, 3).Value = "Ñuble"
$ws.Cells.Item(199, 4).Value = 44246
$ws.Cells.Item(199, 5).Value = 16
$ws.Cells.Item(199, 6).Value = 100114001
$ws.Cells.Item(199, 7).Value = "Papa"
$ws.Cells.Item(199, 8).Value = "Patagonia"
$ws.Cells.Item(199, 9).Value = "1a nueva(o)"
$ws.Cells.Item(199, 10).Value = 120
$ws.Cells.Item(199, 11).Value = 8000
$ws.Cells.Item(199, 12).Value = 8500
$ws.Cells.Item(199, 13).Value = 8250
$ws.Cells.Item(199, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(199, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(199, 16).Value = 330
$ws.Cells.Item(199, 17).Value = 25
$ws.Cells.Item(199, 18).Value = "Hortaliza"

# Update row 198 with this week's new figures (variety/quality/volume/unit
# stay the same; date, prices, origin and $/Kg change).
$ws.Cells.Item(198, 4).Value = 44522
$ws.Cells.Item(198, 11).Value = 9000
$ws.Cells.Item(198, 12).Value = 9500
$ws.Cells.Item(198, 13).Value = 9250
$ws.Cells.Item(198, 15).Value = "Región del Maule"
$ws.Cells.Item(198, 16).Value = 370
